$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear substitute-product rows that are no longer part of the list.
# (These three little groups were removed first, one at a time.)
$ws.Range("C3:D3").ClearContents()
$ws.Range("C4:D4").ClearContents()
$ws.Range("C5:D5").ClearContents()

$ws.Range("C12:D12").ClearContents()
$ws.Range("C13:D13").ClearContents()
$ws.Range("G13").ClearContents()

# Final deletion: the big block of substitutes C16:D31 selected and
# cleared in one go - this also updates the active selection.
$ws.Range("C16:D31").Select()
$ws.Range("C16:D31").ClearContents()

# Rows whose height was driven by the now-removed wrapped text revert to
# the sheet's default (auto) row height.
$rowsToAutoFit = @(18, 20, 21, 22, 23, 25, 26, 27, 28, 29, 30, 31)
foreach ($r in $rowsToAutoFit) {
    $ws.Rows.Item($r).AutoFit()
}
